$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.981.10"
$ws.Range("E2").Value = "  +4.57%  "

$ws.Range("D3").Value = "2.286.72"
$ws.Range("E3").Value = "  +5.13%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'252.33"
$ws.Range("E5").Value = "  +0.27%  "

$ws.Range("D6").Value = "'0.642"
$ws.Range("E6").Value = "  +4.87%  "

$ws.Range("D7").Value = "'72.93"
$ws.Range("E7").Value = "  +10.44%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "'0.663"
$ws.Range("E9").Value = "  +15.30%  "

$ws.Range("D10").Value = "'39.24"
$ws.Range("E10").Value = "  +7.76%  "

$ws.Range("D11").Value = "'0.0982"
$ws.Range("E11").Value = "  +5.20%  "

$ws.Range("D12").Value = "'60.04"
$ws.Range("E12").Value = "  +1.75%  "

$ws.Range("D13").Value = "'7.34"
$ws.Range("E13").Value = "  +7.48%  "

$ws.Range("E14").Value = "  +2.24%  "

$ws.Range("D15").Value = "2.627.01"
$ws.Range("E15").Value = "  +5.18%  "

$ws.Range("D16").Value = "'15.13"
$ws.Range("E16").Value = "  +6.20%  "

$ws.Range("D17").Value = "'0.895"
$ws.Range("E17").Value = "  +6.02%  "

$ws.Range("D18").Value = "2.278.30"
$ws.Range("E18").Value = "  +5.14%  "

$ws.Range("D19").Value = "42.892.44"
$ws.Range("E19").Value = "  +4.49%  "

$ws.Range("D20").Value = "'0.0000102"
$ws.Range("E20").Value = "  +7.46%  "

$ws.Range("E21").Value = "  +5.88%  "

$ws.Range("D22").Value = "'73.72"
$ws.Range("E22").Value = "  +3.11%  "

$ws.Range("D23").Value = "'237.98"
$ws.Range("E23").Value = "  +3.49%  "

$ws.Range("D24").Value = "'2.16"
$ws.Range("E24").Value = "  +7.64%  "

$ws.Range("E25").Value = "  +1.73%  "

$ws.Range("D26").Value = "'11.71"
$ws.Range("E26").Value = "  +3.41%  "

$ws.Range("E27").Value = "  +0.09%  "

$ws.Range("E28").Value = "  +2.27%  "

$ws.Range("E29").Value = "  -0.78%  "

$ws.Range("E30").Value = "  +0.85%  "

$ws.Range("D31").Value = "'168.53"
$ws.Range("E31").Value = "  +0.40%  "

$ws.Range("D32").Value = "'21.18"
$ws.Range("E32").Value = "  +5.17%  "

$ws.Range("D33").Value = "'6.32"
$ws.Range("E33").Value = "  +11.51%  "

$ws.Range("D34").Value = "'0.129"
$ws.Range("E34").Value = "  +7.26%  "

$ws.Range("D35").Value = "'0.0816"
$ws.Range("E35").Value = "  +9.23%  "

$ws.Range("D36").Value = "'31.52"
$ws.Range("E36").Value = "  +28.90%  "

$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D37").Value = "'0.127"
$ws.Range("E37").Value = "  +5.22%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'4.80"
$ws.Range("E38").Value = "  +22.24%  "

$ws.Range("E39").Value = "  +6.94%  "

$ws.Range("D40").Value = "'0.0312"
$ws.Range("E40").Value = "  +2.88%  "

$ws.Range("E41").Value = "  +5.73%  "

$ws.Range("D42").Value = "'13.25"
$ws.Range("E42").Value = "  +17.06%  "

$ws.Range("D43").Value = "'6.07"
$ws.Range("E43").Value = "  +11.00%  "

$ws.Range("D44").Value = "'0.214"
$ws.Range("E44").Value = "  +14.36%  "

$ws.Range("E45").Value = "  +9.18%  "

$ws.Range("E46").Value = "  -9.30%  "

$ws.Range("D47").Value = "'61.97"
$ws.Range("E47").Value = "  +1.73%  "

$ws.Range("E48").Value = "  +5.34%  "

$ws.Range("E49").Value = "  +5.66%  "

$ws.Range("E50").Value = "  +0.11%  "

$ws.Range("E51").Value = "  +5.39%  "
